# Add a new "Spain" market sheet, cloned from the existing "Italy" sheet,
# with the Spain-specific market name / NGC code filled in (mirrors how
# every other country tab in this workbook was produced).

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Duplicate the Italy sheet (placed right after it) - this carries over all
# styles, merged cells, column widths, etc. exactly like the other tabs.
$italy.Copy([System.Reflection.Missing]::Value, $italy) | Out-Null
$spain = $wb.ActiveSheet
$spain.Name = "Spain"

# Fill in the Spain-specific data.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064/T2063/T2062/T2065/T2056"

# The new, longer text in column B (and the resulting wrap in column D)
# means the columns/rows no longer fit their old Italy-sized dimensions -
# resize them like Excel's own AutoFit would.
$spain.Columns("B").ColumnWidth = 47.41666666666731
$spain.Columns("C").ColumnWidth = 14.916666666666885
$spain.Columns("D").ColumnWidth = 23.25000000000003
$spain.Rows("3:5").RowHeight = 28.8

# Restore a sane selection on the Italy sheet (it had been left on the
# blank row below the table) and leave the new Spain sheet active with its
# own selection, matching the end state of the edit.
$italy.Range("A1:D11").Select() | Out-Null
$spain.Select() | Out-Null
$spain.Range("B10").Select() | Out-Null
